$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that will receive numeric-looking text values need the cell
# pre-formatted as Text so Excel stores them as strings (preserving
# leading/trailing zeros etc.) instead of coercing to a Number.
$textCells = @("D5","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17","D24","D25","D26","D29","D30","D31","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '30.110.13'
$ws.Range("E2").Value = '  -0.55%  '

# Row 3
$ws.Range("D3").Value = '1.856.40'
$ws.Range("E3").Value = '  -0.64%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").Value = '233.24'
$ws.Range("E5").Value = '  -0.85%  '

# Row 6
$ws.Range("E6").Value = '  +0.02%  '

# Row 7
$ws.Range("D7").Value = '0.4690'
$ws.Range("E7").Value = '  -0.14%  '

# Row 8
$ws.Range("D8").Value = '42.87'
$ws.Range("E8").Value = '  -0.09%  '

# Row 9
$ws.Range("D9").Value = '0.2822'

# Row 10
$ws.Range("D10").Value = '0.06443'
$ws.Range("E10").Value = '  -2.03%  '

# Row 11
$ws.Range("D11").Value = '20.92'
$ws.Range("E11").Value = '  -4.05%  '

# Row 12
$ws.Range("D12").Value = '0.07733'

# Row 13
$ws.Range("D13").Value = '1.859.72'
$ws.Range("E13").Value = '  -0.55%  '

# Row 14
$ws.Range("D14").Value = '93.25'
$ws.Range("E14").Value = '  -4.00%  '

# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '5.034'
$ws.Range("E15").Value = '  -1.74%  '

# Row 16
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = '0.6762'
$ws.Range("E16").Value = '  -1.27%  '

# Row 17
$ws.Range("D17").Value = '265.70'
$ws.Range("E17").Value = '  -1.09%  '

# Row 18
$ws.Range("D18").Value = '30.086.97'
$ws.Range("E18").Value = '  -0.59%  '

# Row 19
$ws.Range("E19").Value = '  -5.41%  '

# Row 20
$ws.Range("E20").Value = '  -1.49%  '

# Row 22
$ws.Range("D22").Value = '2.112.84'
$ws.Range("E22").Value = '  -0.15%  '

# Row 23
$ws.Range("E23").Value = '  +0.07%  '

# Row 24
$ws.Range("D24").Value = '5.124'

# Row 25
$ws.Range("D25").Value = '6.083'
$ws.Range("E25").Value = '  -2.06%  '

# Row 26
$ws.Range("D26").Value = '9.265'
$ws.Range("E26").Value = '  -1.61%  '

# Row 27
$ws.Range("E27").Value = '  -1.77%  '

# Row 28
$ws.Range("E28").Value = '  -2.34%  '

# Row 29
$ws.Range("D29").Value = '1.875'
$ws.Range("E29").Value = '  -3.84%  '

# Row 30
$ws.Range("D30").Value = '1.364'
$ws.Range("E30").Value = '  -0.35%  '

# Row 31
$ws.Range("D31").Value = '0.09809'
$ws.Range("E31").Value = '  -0.55%  '

# Row 32
$ws.Range("E32").Value = '  -0.90%  '

# Row 33
$ws.Range("D33").Value = '4.182'
$ws.Range("E33").Value = '  -4.50%  '

# Row 34
$ws.Range("D34").Value = '3.958'
$ws.Range("E34").Value = '  -2.87%  '

# Row 35
$ws.Range("D35").Value = '0.04627'
$ws.Range("E35").Value = '  -1.71%  '

# Row 36
$ws.Range("D36").Value = '1.109'
$ws.Range("E36").Value = '  -2.07%  '

# Row 37
$ws.Range("D37").Value = '0.6844'
$ws.Range("E37").Value = '  -2.30%  '

# Row 38
$ws.Range("D38").Value = '2.713'
$ws.Range("E38").Value = '  +0.08%  '

# Row 39
$ws.Range("D39").Value = '0.01822'
$ws.Range("E39").Value = '  -2.76%  '

# Row 40
$ws.Range("D40").Value = '2.708'
$ws.Range("E40").Value = '  +3.18%  '

# Row 41
$ws.Range("D41").Value = '6.231'
$ws.Range("E41").Value = '  -1.04%  '

# Row 42
$ws.Range("D42").Value = '70.35'
$ws.Range("E42").Value = '  -2.66%  '

# Row 43
$ws.Range("E43").Value = '  +0.03%  '

# Row 44
$ws.Range("D44").Value = '0.8304'
$ws.Range("E44").Value = '  -1.42%  '

# Row 45
$ws.Range("D45").Value = '101.66'
$ws.Range("E45").Value = '  -1.36%  '

# Row 46
$ws.Range("D46").Value = '1.858'
$ws.Range("E46").Value = '  -4.78%  '

# Row 47
$ws.Range("D47").Value = '0.4025'

# Row 48
$ws.Range("D48").Value = '9.121'
$ws.Range("E48").Value = '  -0.90%  '

# Row 49
$ws.Range("D49").Value = '6.904'
$ws.Range("E49").Value = '  -2.19%  '

# Row 50
$ws.Range("D50").Value = '917.26'
$ws.Range("E50").Value = '  -0.78%  '

# Row 51
$ws.Range("D51").Value = '33.93'
$ws.Range("E51").Value = '  -1.68%  '

# Restore default style on the text-forced cells (keeps the value as text
# while dropping the now-unused custom number format from those cells).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
